# Estadisticos Matutinos 15 Oct
# Update D:K statistics for "1er Parcial" and "3er Parcial" sheets (rows 2-12).

$wb = $excel.ActiveWorkbook

# row -> E, F, G, H, I, J, K
$data = @{
    2  = @(9,  15, 37.5,               62.5,  7.8, 15, 62.5)
    3  = @(21, 11, 65.63,              34.38, 7.8, 11, 34.38)
    4  = @(14, 20, 41.18,              58.82, 8.4, 20, 58.82)
    5  = @(8,  13, 38.1,               61.9,  7.9, 13, 61.9)
    6  = @(26, 9,  74.29000000000001,  25.71, 7.8, 9,  25.71)
    7  = @(8,  13, 38.1,               61.9,  7.4, 13, 61.9)
    8  = @(21, 10, 67.73999999999999,  32.26, 7.1, 0,  0)
    9  = @(18, 17, 51.43,              48.57, 6.8, 17, 48.57)
    10 = @(34, 5,  87.18000000000001,  12.82, 7.7, 0,  0)
    11 = @(30, 5,  85.70999999999999,  14.29, 8.1, 5,  14.29)
    12 = @(16, 17, 48.48,              51.52, 7.2, 17, 51.52)
}

foreach ($sheetName in @("1er Parcial", "3er Parcial")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $data.Keys) {
        $vals = $data[$row]
        $ws.Cells.Item($row, 5).Value  = $vals[0]   # E: Aprobados
        $ws.Cells.Item($row, 6).Value  = $vals[1]   # F: Reprobados
        $ws.Cells.Item($row, 7).Value  = $vals[2]   # G: Por_Apro
        $ws.Cells.Item($row, 8).Value  = $vals[3]   # H: Por_Repro
        $ws.Cells.Item($row, 9).Value  = $vals[4]   # I: Promedio
        $ws.Cells.Item($row, 10).Value = $vals[5]   # J: Blancos
        $ws.Cells.Item($row, 11).Value = $vals[6]   # K: Por_Blan
    }
}
